$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff shows a cyclic rotation of the D, M, N, O, P, S column values
# across rows 2-5 (row data for Fecha/Volumen/Precio min/max/promedio/Precio $/Kg):
#   new row2 = old row5
#   new row3 = old row4
#   new row4 = old row2
#   new row5 = old row3
# Capture original values first, then write them back in rotated order.

$cols = @("D","M","N","O","P","S")

$orig = @{}
foreach ($r in 2..5) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

$mapping = @{ 2 = 5; 3 = 4; 4 = 2; 5 = 3 }

foreach ($r in 2..5) {
    $src = $mapping[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $orig[$src][$c]
    }
}
